$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K shifts to F:M)
$ws.Columns("D:E").Insert()

# Copy cell formatting (number formats/styles) from column F into new D:E columns,
# processing in separate row-blocks to avoid touching rows that have no cells at all
# in that block (title rows 5/6/37/79 and the blank separator rows 36/78)
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns (most recent two quarters) and correct the
# restated figures that moved into column F (previously column D).
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 113000
$ws.Range("E8").Value = 118400
$ws.Range("F8").Value = 251600
$ws.Range("D9").Value = 127300
$ws.Range("E9").Value = 98200
$ws.Range("F9").Value = 193500
$ws.Range("D10").Value = -14300
$ws.Range("E10").Value = 20200
$ws.Range("F10").Value = 58100
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 6400
$ws.Range("E14").Value = 2300
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 147600
$ws.Range("E17").Value = 114100
$ws.Range("F17").Value = 225300
$ws.Range("D18").Value = -34600
$ws.Range("E18").Value = 4300
$ws.Range("F18").Value = 26300
$ws.Range("D20").Value = 57700
$ws.Range("E20").Value = -5300
$ws.Range("F20").Value = 23600
$ws.Range("D21").Value = 115400
$ws.Range("E21").Value = 76300
$ws.Range("F21").Value = 160400
$ws.Range("D22").Value = 28400
$ws.Range("E22").Value = 27500
$ws.Range("F22").Value = 53200
$ws.Range("D23").Value = -5300
$ws.Range("E23").Value = -28500
$ws.Range("F23").Value = -3200
$ws.Range("D24").Value = 17800
$ws.Range("E24").Value = 3000
$ws.Range("F24").Value = 11200
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -23100
$ws.Range("E26").Value = -31500
$ws.Range("F26").Value = -14400
$ws.Range("D27").Value = -13900
$ws.Range("E27").Value = -12600
$ws.Range("F27").Value = 168500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -57700
$ws.Range("E32").Value = 5300
$ws.Range("F32").Value = -23600
$ws.Range("D33").Value = -13900
$ws.Range("E33").Value = -12600
$ws.Range("F33").Value = 168500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -13900
$ws.Range("E35").Value = -12600
$ws.Range("F35").Value = 168500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 101000
$ws.Range("E41").Value = 125700
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 50000
$ws.Range("E43").Value = 50300
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 60000
$ws.Range("E45").Value = 64200
$ws.Range("D46").Value = 211000
$ws.Range("E46").Value = 240200
$ws.Range("D47").Value = 270000
$ws.Range("E47").Value = 372400
$ws.Range("D48").Value = 4378000
$ws.Range("E48").Value = 4331000
$ws.Range("D49").Value = 277000
$ws.Range("E49").Value = 277600
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 158000
$ws.Range("E52").Value = 150800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 5294000
$ws.Range("E54").Value = 5372100
$ws.Range("D57").Value = 67000
$ws.Range("E57").Value = 57400
$ws.Range("D58").Value = 254000
$ws.Range("E58").Value = 250000
$ws.Range("D59").Value = 157000
$ws.Range("E59").Value = 134000
$ws.Range("D60").Value = 478000
$ws.Range("E60").Value = 441400
$ws.Range("D61").Value = 2029000
$ws.Range("E61").Value = 2129600
$ws.Range("D62").Value = 628000
$ws.Range("E62").Value = 616500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 4247000
$ws.Range("E66").Value = 4232900
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -27000
$ws.Range("E72").Value = -12600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1047000
$ws.Range("E76").Value = 1139200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -13900
$ws.Range("E81").Value = -12600
$ws.Range("F81").Value = 168500
$ws.Range("D83").Value = 92300
$ws.Range("E83").Value = 77300
$ws.Range("F83").Value = 110400
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 48500
$ws.Range("E89").Value = 120500
$ws.Range("F89").Value = 109900
$ws.Range("D91").Value = -170800
$ws.Range("E91").Value = -2800
$ws.Range("F91").Value = -7400
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -29400
$ws.Range("E94").Value = -47200
$ws.Range("F94").Value = -296400
$ws.Range("D96").Value = -41400
$ws.Range("E96").Value = -41100
$ws.Range("F96").Value = -82500
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -41100
$ws.Range("E100").Value = -57700
$ws.Range("F100").Value = 181800
$ws.Range("D101").Value = -1000
$ws.Range("E101").Value = -600
$ws.Range("F101").Value = -2400
$ws.Range("D102").Value = -22900
$ws.Range("E102").Value = 15000
$ws.Range("F102").Value = -7100

Write-Host "Edit complete"
